$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Library_Formula")

# Add a new row (row 13) mirroring row 12, with new action/description text.
$ws.Range("A13").Value = "CREATE/MODIFY"
$ws.Range("B13").Value = "WeightedModelWoe_PageLib"
$ws.Range("C13").Value = "updateVariables"
$ws.Range("E13").Value = "String"
$ws.Range("F13").Value = "String, String"

$ws.Range("A13").Style = $ws.Range("A12").Style
$ws.Range("B13").Style = $ws.Range("B12").Style
$ws.Range("C13").Style = $ws.Range("C12").Style
$ws.Range("E13").Style = $ws.Range("E12").Style
$ws.Range("F13").Style = $ws.Range("F12").Style

$ws.Range("F14").Select()
